$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.075.55"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "1.680.38"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("E6").Value = "  -3.26%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +1.84%  "
$ws.Range("E9").Value = "  +5.54%  "
$ws.Range("E10").Value = "  +0.56%  "
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").Value = "1.917.33"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").Value = "1.711.94"
$ws.Range("E13").Value = "  +2.71%  "
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.39"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").Value = "27.062.50"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("E18").Value = "  +2.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "236.54"
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("E22").Value = "  +2.10%  "
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("E24").Value = "  -2.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.05"
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("E26").Value = "  +2.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.51"
$ws.Range("E27").Value = "  +3.83%  "
$ws.Range("E28").Value = "  -1.61%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").Value = "1.545.34"
$ws.Range("E33").Value = "  +6.27%  "
$ws.Range("E34").Value = "  +1.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.73"
$ws.Range("E35").Value = "  +4.63%  "
$ws.Range("E36").Value = "  -1.15%  "
$ws.Range("E37").Value = "  +1.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.917"
$ws.Range("E38").Value = "  +1.37%  "
$ws.Range("E39").Value = "  +2.37%  "
$ws.Range("E40").Value = "  +7.07%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "67.83"
$ws.Range("E42").Value = "  +2.52%  "
$ws.Range("E43").Value = "  -3.28%  "
$ws.Range("E44").Value = "  -0.90%  "
$ws.Range("D45").Value = "1.822.16"
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.50"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("E48").Value = "  +2.91%  "
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("E50").Value = "  +1.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.06"
$ws.Range("E51").Value = "  +7.02%  "
